$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 15 (DRA052) with the new description and Jira id.
# Set C15 (description) before B15 (Jira id) so new shared strings are
# appended in the same order as in the target workbook (description=45, jira id=46).
$ws.Range("C15").Value = "Verify that, the account link or merge should not be made, if skipping the linking of an existing social account"
$ws.Range("B15").Value = "OPQA-4243"

# Row 15 height shrinks from 60 to 30 now that the description text is shorter.
$ws.Rows.Item(15).RowHeight = 30

# Update the active selection to C15.
$ws.Range("C15").Select()
